# Atualiza planilha de análise de curto prazo da Embraer
#
# 1) "Indicadores Principais": recompute Liquidez_Seca (col C) to account for
#    a newly-recognised "DespesasAntecipadas" deduction.
# 2) "Todos Indicadores": same Liquidez Seca recompute (col C), plus a tiny
#    floating point refresh of CCL (PL+PNC-ANC) (col F) caused by the
#    AtivoNaoCirculante recompute on "Dados Base".
# 3) "Dados Base": restructure columns - move AtivoNaoCirculante (was col T)
#    to right after AtivoCirculante (new col G), and insert a brand new
#    DespesasAntecipadas column (new col I) with newly supplied figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Dados Base": column restructuring + new data
# ---------------------------------------------------------------------------
$wsBase = $wb.Worksheets.Item("Dados Base")

# Move column T (AtivoNaoCirculante) so it sits right after AtivoCirculante
# (current column F), i.e. becomes the new column G. This shifts the old
# G:S block one column to the right (new H:T).
$wsBase.Columns("T:T").Cut()
$wsBase.Columns("G:G").Insert()

# Insert a brand new blank column at I (right after the relocated
# AtivoTotal, now column H) to hold DespesasAntecipadas. This shifts the
# old H:T block (EmprestimosFinanciamentosCP .. PassivosCiclicosCP) one more
# column to the right, landing at J:U.
$wsBase.Columns("I:I").Insert()

$wsBase.Range("I1").Value = "DespesasAntecipadas"

$despesasAntecipadas = @{ 2 = 205377.4376749034; 3 = 212494.6416056913; 4 = 342683.40968184; 5 = 252324.6084; 6 = 363721 }
$ativoNaoCirculante  = @{ 2 = 29141300.63507662; 3 = 27650888.97605425; 4 = 25044728.79473388; 5 = 23942405.2206; 6 = 32892812 }

foreach ($row in 2..6) {
    $wsBase.Cells.Item($row, 7).Value = $ativoNaoCirculante[$row]   # col G = AtivoNaoCirculante
    $wsBase.Cells.Item($row, 9).Value = $despesasAntecipadas[$row]  # col I = DespesasAntecipadas
}

# ---------------------------------------------------------------------------
# Sheet "Indicadores Principais": refresh Liquidez_Seca (col C)
# ---------------------------------------------------------------------------
$wsPrincipais = $wb.Worksheets.Item("Indicadores Principais")

$liquidezSeca = @{ 2 = 1.501536121434408; 3 = 1.363416501871598; 4 = 1.04674837478906; 5 = 0.9192568588539315; 6 = 0.796505439569714 }

foreach ($row in 2..6) {
    $wsPrincipais.Cells.Item($row, 3).Value = $liquidezSeca[$row]  # col C = Liquidez_Seca
}

# ---------------------------------------------------------------------------
# Sheet "Todos Indicadores": refresh Liquidez Seca (col C) and
# CCL (PL+PNC-ANC) (col F)
# ---------------------------------------------------------------------------
$wsTodos = $wb.Worksheets.Item("Todos Indicadores")

foreach ($row in 2..6) {
    $wsTodos.Cells.Item($row, 3).Value = $liquidezSeca[$row]  # col C = Liquidez Seca
}

$cclPlPncAnc = @{ 2 = 24465118.76410241; 3 = 19688354.08991012; 4 = 14502530.98141576 }
foreach ($row in 2..4) {
    $wsTodos.Cells.Item($row, 6).Value = $cclPlPncAnc[$row]  # col F = CCL (PL+PNC-ANC)
}
